$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing "3 + 2" / "28-10-23 & 29-10-23" entries (row 33, F/G) ---
$ws.Range("F33").Value = "3 + 2 + 1.5"
$ws.Range("G33").Value = "28-10-23 to 30-10-23"

# --- Fill in Expected WH / Actual WH / Date for rows 35-42 (SQL training days) ---
# Row 35
$ws.Range("E35").Value = 1.5
$ws.Range("F35").Value = 1
$ws.Range("G35").Value = 44937
$ws.Range("G17").Copy()
$ws.Range("G35").PasteSpecial(-4122)

# Row 36
$ws.Range("E36").Value = 1
$ws.Range("F36").Value = 0.75
$ws.Range("G36").Value = 44937
$ws.Range("G17").Copy()
$ws.Range("G36").PasteSpecial(-4122)

# Row 37
$ws.Range("E37").Value = 1.5
$ws.Range("F37").Value = 1
$ws.Range("G37").Value = 44968
$ws.Range("G17").Copy()
$ws.Range("G37").PasteSpecial(-4122)

# Row 38
$ws.Range("E38").Value = 2
$ws.Range("F38").Value = 1
$ws.Range("G38").Value = 45057
$ws.Range("G17").Copy()
$ws.Range("G38").PasteSpecial(-4122)

# Row 39
$ws.Range("E39").Value = 2
$ws.Range("F39").Value = 1.5
$ws.Range("G39").Value = 44968
$ws.Range("G17").Copy()
$ws.Range("G39").PasteSpecial(-4122)

# Row 40 (date is stored as literal text, not a real date serial)
$ws.Range("E40").Value = 1
$ws.Range("F40").Value = 0.5
$ws.Range("G40").Value = "3/11/023"

# Row 41
$ws.Range("E41").Value = 1
$ws.Range("F41").Value = 0.5
$ws.Range("G41").Value = 44996
$ws.Range("G17").Copy()
$ws.Range("G41").PasteSpecial(-4122)

# Row 42
$ws.Range("E42").Value = 1
$ws.Range("F42").Value = 1
$ws.Range("G42").Value = 45027
$ws.Range("G17").Copy()
$ws.Range("G42").PasteSpecial(-4122)

# --- Rows 43-51: only "Expected WH" (column E) filled in ---
$ws.Range("E43").Value = 2.5
$ws.Range("E44").Value = 2
$ws.Range("E45").Value = 1.5
$ws.Range("E46").Value = 2
$ws.Range("E47").Value = 1.5
$ws.Range("E48").Value = 1.5
$ws.Range("E49").Value = 3
$ws.Range("E50").Value = 3
$ws.Range("E51").Value = 3

# --- Update view state: scrolled down to row 33, selection on G38 ---
$excel.ActiveWindow.ScrollRow = 33
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G38").Select()
